# Regenerate the localization-status handoff report.
#
# The underlying CI job re-ran and produced a fresh snapshot: the
# "334ba757-..." file moved from "In Translation" to "Ready for handoff"
# (new handoff timestamps / new xlf names), while the "88989a8a-..." file's
# status stayed the same. The two files also swapped which row they land
# on in every sheet. We replicate this by rewriting each worksheet's data
# rows and re-pointing the hyperlinks (delete + re-add, so every row keeps
# exactly one hyperlink per hyperlinked cell, matching the file it now
# displays).

$wb = $excel.ActiveWorkbook

$mdBase   = "https://github.com/OpenLocalizationTest/oltest/blob/e7caed04040e78178c84fde070f8fa349131ec5f/e2e/"
$zhHtBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/19472b3e8eca5dc664fcf83abefbb090c939b1b0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/"
$deHtBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ff3c6595268786a200132faa91ccd69cc752491b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/"

$file1 = "334ba757-3e04-4960-b90c-ec9a3eb937d1"
$file2 = "88989a8a-0262-4447-ab47-0b581d8ecc69"

$file1MdUrl = $mdBase + $file1 + ".md"
$file2MdUrl = $mdBase + $file2 + ".md"

$file1ZhXlf = $file1 + ".f1c51fedfdf63615d256bce43d1c1a62ece4c6bb.zh-cn.xlf"
$file2ZhXlf = $file2 + ".8161ab9136248d51b8caa3fafd358546e9cb76a5.zh-cn.xlf"
$file1DeXlf = $file1 + ".f1c51fedfdf63615d256bce43d1c1a62ece4c6bb.de-de.xlf"
$file2DeXlf = $file2 + ".8161ab9136248d51b8caa3fafd358546e9cb76a5.de-de.xlf"

$file1ZhXlfUrl = $zhHtBase + $file1ZhXlf
$file2ZhXlfUrl = $zhHtBase + $file2ZhXlf
$file1DeXlfUrl = $deHtBase + $file1DeXlf
$file2DeXlfUrl = $deHtBase + $file2DeXlf

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Hyperlinks.Delete()

$wsOverview.Range("A2").Value = $file2 + ".md"
$wsOverview.Range("B2").Value = "In Translation"
$wsOverview.Range("C2").Value = "In Translation"
$wsOverview.Range("D2").Value = "2016-12-11 22:12:26"

$wsOverview.Range("A3").Value = $file1 + ".md"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-13-11 22:13:08"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $file2MdUrl, "", "", $file2 + ".md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $file1MdUrl, "", "", $file1 + ".md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | File Extension | Status | Latest
# Handoff File | Latest Handoff Datetime | ... | Dependency From | Error Detail
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Hyperlinks.Delete()

$wsZh.Range("A2").Value = $file2 + ".md"
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = "In Translation"
$wsZh.Range("D2").Value = $file2ZhXlf
$wsZh.Range("E2").Value = "2016-03-11 22:12:16"
$wsZh.Range("H2").Value = "0001-01-01 00:00:00"
$wsZh.Range("I2").Value = "Include"

$wsZh.Range("A3").Value = $file1 + ".md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = $file1ZhXlf
$wsZh.Range("E3").Value = "2016-03-11 22:13:05"
$wsZh.Range("H3").Value = "0001-01-01 00:00:00"
$wsZh.Range("I3").Value = "Include"

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $file2MdUrl, "", "", $file2 + ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), $file2MdUrl, "", "", ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $file2ZhXlfUrl, "", "", $file2ZhXlf) | Out-Null

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $file1MdUrl, "", "", $file1 + ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B3"), $file1MdUrl, "", "", ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $file1ZhXlfUrl, "", "", $file1ZhXlf) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de": same column layout as zh-cn
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Hyperlinks.Delete()

$wsDe.Range("A2").Value = $file2 + ".md"
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = "In Translation"
$wsDe.Range("D2").Value = $file2DeXlf
$wsDe.Range("E2").Value = "2016-03-11 22:12:26"
$wsDe.Range("H2").Value = "0001-01-01 00:00:00"
$wsDe.Range("I2").Value = "Include"

$wsDe.Range("A3").Value = $file1 + ".md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = $file1DeXlf
$wsDe.Range("E3").Value = "2016-03-11 22:13:08"
$wsDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDe.Range("I3").Value = "Include"

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $file2MdUrl, "", "", $file2 + ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), $file2MdUrl, "", "", ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $file2DeXlfUrl, "", "", $file2DeXlf) | Out-Null

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $file1MdUrl, "", "", $file1 + ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B3"), $file1MdUrl, "", "", ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $file1DeXlfUrl, "", "", $file1DeXlf) | Out-Null

Write-Output "Report regenerated for handoff."
